# "Generate Report for Handoff"
#
# The localization status report is regenerated: the zh-cn / de-de rows
# move from "In Translation" to "Ready for handoff", the associated
# timestamps advance to the new generation time, and the Status /
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" columns
# widen to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status column
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status column
$wsZhCn.Range("C2").Value = "Ready for handoff"       # Status column
$wsDeDe.Range("C2").Value = "Ready for handoff"       # Status column

# --- Timestamps bump forward with the new handoff generation -----------
$wsOverview.Range("G2").Value = "2016-09-05 15:13:21" # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value    = "2016-09-05 15:13:21"  # Latest Handoff Datetime (de-de)
$wsZhCn.Range("H2").Value   = "2016-09-05 15:13:17"   # Latest Handoff Datetime (zh-cn)

# --- Widen the columns that now hold the longer status/date text -------
$newWidth = 17.2159881591797 - 0.8333333333333334

$wsOverview.Range("E1").ColumnWidth = $newWidth
$wsOverview.Range("F1").ColumnWidth = $newWidth
$wsZhCn.Range("C1").ColumnWidth = $newWidth
$wsDeDe.Range("C1").ColumnWidth = $newWidth
